# get toi_source and API data for ncp-gop from eims-toi output
#
# Insert a new "inorganic matter" keyword row into the Keywords sheet
# (pushing the existing chemistry/oceanography/seawater/dissolved oxygen
# rows down by one), and make the Keywords sheet the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keywords")

# Insert a new row above the current row 3 ("chemistry"), shifting the
# existing keyword rows (and the two trailing blank rows) down by one.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the new keyword entry.
$ws.Cells.Item(3, 1).Value = "inorganic matter"
$ws.Cells.Item(3, 2).Value = "LTER Core Research Areas"

# Make the Keywords sheet the active sheet/tab, with A3:B3 selected.
$ws.Activate()
$ws.Range("A3:B3").Select()
